$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Glg1"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 21.16594
$ws.Range("H2").Value = 63.49782
$ws.Range("I2").Value = 0.1484747310246988
$ws.Range("J2").Value = 0.1568390841279916
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.112632333333332
$ws.Range("N2").Value = 27.337897
$ws.Range("O2").Value = 0.9981738658344552
$ws.Range("P2").Value = 0.9981738658344552
$ws.Range("Q2").Value = 192.8774292093933
$ws.Range("R2").Value = 1735.89686288454
$ws.Range("S2").Value = 0.1482035962456545
$ws.Range("T2").Value = 0.1565526749179727

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Glg1"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 21.16594
$ws.Range("H3").Value = 63.49782
$ws.Range("I3").Value = 0.1484747310246988
$ws.Range("J3").Value = 0.1568390841279916
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01667133333333333
$ws.Range("N3").Value = 0.050014
$ws.Range("O3").Value = 0.001826134165544791
$ws.Range("P3").Value = 0.001826134165544791
$ws.Range("Q3").Value = 0.3528644410533333
$ws.Range("R3").Value = 3.17577996948
$ws.Range("S3").Value = 0.0002711347790442757
$ws.Range("T3").Value = 0.0002864092100188792

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Glg1"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 44.240478
$ws.Range("H4").Value = 132.721434
$ws.Range("I4").Value = 0.3103378858417866
$ws.Range("J4").Value = 0.3278208315295499
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.112632333333332
$ws.Range("N4").Value = 27.337897
$ws.Range("O4").Value = 0.9981738658344552
$ws.Range("P4").Value = 0.9981738658344552
$ws.Range("Q4").Value = 403.1472102649219
$ws.Range("R4").Value = 3628.324892384298
$ws.Range("S4").Value = 0.309771167225588
$ws.Range("T4").Value = 0.3272221867089165

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Glg1"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 44.240478
$ws.Range("H5").Value = 132.721434
$ws.Range("I5").Value = 0.3103378858417866
$ws.Range("J5").Value = 0.3278208315295499
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01667133333333333
$ws.Range("N5").Value = 0.050014
$ws.Range("O5").Value = 0.001826134165544791
$ws.Range("P5").Value = 0.001826134165544791
$ws.Range("Q5").Value = 0.7375477555639999
$ws.Range("R5").Value = 6.637929800076
$ws.Range("S5").Value = 0.0005667186161986257
$ws.Range("T5").Value = 0.0005986448206334142

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Glg1"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.368218
$ws.Range("H6").Value = 70.104654
$ws.Range("I6").Value = 0.1639232598257637
$ws.Range("J6").Value = 0.1731579088300944
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.112632333333332
$ws.Range("N6").Value = 27.337897
$ws.Range("O6").Value = 0.9981738658344552
$ws.Range("P6").Value = 0.9981738658344552
$ws.Range("Q6").Value = 212.945978919182
$ws.Range("R6").Value = 1916.513810272638
$ws.Range("S6").Value = 0.1636239139604684
$ws.Range("T6").Value = 0.1728416992567454

$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Glg1"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 23.368218
$ws.Range("H7").Value = 70.104654
$ws.Range("I7").Value = 0.1639232598257637
$ws.Range("J7").Value = 0.1731579088300944
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01667133333333333
$ws.Range("N7").Value = 0.050014
$ws.Range("O7").Value = 0.001826134165544791
$ws.Range("P7").Value = 0.001826134165544791
$ws.Range("Q7").Value = 0.389579351684
$ws.Range("R7").Value = 3.506214165156
$ws.Range("S7").Value = 0.000299345865295303
$ws.Range("T7").Value = 0.0003162095733489254

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Glg1"
$ws.Range("C8").Value = "Sele"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 30.973355
$ws.Range("H8").Value = 92.92006500000001
$ws.Range("I8").Value = 0.2172717371662922
$ws.Range("J8").Value = 0.2295117831086713
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.112632333333332
$ws.Range("N8").Value = 27.337897
$ws.Range("O8").Value = 0.9981738658344552
$ws.Range("P8").Value = 0.9981738658344552
$ws.Range("Q8").Value = 282.2487962448116
$ws.Range("R8").Value = 2540.239166203305
$ws.Range("S8").Value = 0.2168749698238456
$ws.Range("T8").Value = 0.2290926638001414

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Glg1"
$ws.Range("C9").Value = "Sele"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 30.973355
$ws.Range("H9").Value = 92.92006500000001
$ws.Range("I9").Value = 0.2172717371662922
$ws.Range("J9").Value = 0.2295117831086713
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01667133333333333
$ws.Range("N9").Value = 0.050014
$ws.Range("O9").Value = 0.001826134165544791
$ws.Range("P9").Value = 0.001826134165544791
$ws.Range("Q9").Value = 0.5163671256566666
$ws.Range("R9").Value = 4.64730413091
$ws.Range("S9").Value = 0.0003967673424466342
$ws.Range("T9").Value = 0.0004191193085298505

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Glg1"
$ws.Range("C10").Value = "Sele"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 22.807849
$ws.Range("H10").Value = 45.61569799999999
$ws.Range("I10").Value = 0.1599923861414587
$ws.Range("J10").Value = 0.1126703924036929
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.112632333333332
$ws.Range("N10").Value = 27.337897
$ws.Range("O10").Value = 0.9981738658344552
$ws.Range("P10").Value = 0.9981738658344552
$ws.Range("Q10").Value = 207.8395422511843
$ws.Range("R10").Value = 1247.037253507106
$ws.Range("S10").Value = 0.1597002185788987
$ws.Range("T10").Value = 0.1124646411506792

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Glg1"
$ws.Range("C11").Value = "Sele"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 22.807849
$ws.Range("H11").Value = 45.61569799999999
$ws.Range("I11").Value = 0.1599923861414587
$ws.Range("J11").Value = 0.1126703924036929
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01667133333333333
$ws.Range("N11").Value = 0.050014
$ws.Range("O11").Value = 0.001826134165544791
$ws.Range("P11").Value = 0.001826134165544791
$ws.Range("Q11").Value = 0.3802372532953333
$ws.Range("R11").Value = 2.281423519772
$ws.Range("S11").Value = 0.0002921675625599526
$ws.Range("T11").Value = 0.0002057512530137219
